$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 20: "Last Job Date" / "January 1st 2022"
$ws.Range("A20").Value = "Last Job Date"
$ws.Range("B20").Value = "January 1st 2022"

# Match the style of the row above (A column) for the label cell
$ws.Range("A20").Style = $ws.Range("A19").Style

# B20 uses a new style with font + fill + border applied (distinct from B19's hyperlink style)
$ws.Range("B20").Font.Size = 10
$ws.Range("B20").Font.Name = "Arial"

# Move the selection to B20, matching the new active cell in the diff
$ws.Range("B20").Select()
